$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Friday's time entries (row 17): In 09:00, Out 12:00
$ws.Range("C17").Value = 0.375
$ws.Range("D17").Value = 0.5

# D17 previously had no explicit time formatting (it was still blank); give it
# the same "Out" time format/style already used by the cells below it
$ws.Range("D18").Copy()
$ws.Range("D17").PasteSpecial(-4122)

# Regular hours formula for Friday, matching the pattern used by the other rows
$ws.Range("G17").Formula = "=IF((((D17-C17)+(F17-E17))*24)>8,8,((D17-C17)+(F17-E17))*24)"

# Note documenting Friday's split shift
$ws.Range("L17").Value = "0900 – 1000, 1300 – 1500"

# Extend the time-format data validation down through row 19 (now that row 17 is fully used)
$ws.Range("C13:F19").Validation.Delete()
$ws.Range("C13:F19").Validation.Add(5, 1, 1, 0, 0.999305555555556)
$ws.Range("C13:F19").Validation.ErrorTitle = "Invalid Entry"
$ws.Range("C13:F19").Validation.ErrorMessage = "Please enter time in military time format between 0:00 and 23:59 (1:00, 8:00, 13:00, 20:00, etc.)."

# Move the active selection to F18, matching the author's cursor position after the edit
[void]$ws.Range("F18").Select()
